$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh per-row coin price/volume figures from the latest feed snapshot.
# Some Price values keep their original look-like-a-number text form
# (e.g. thousands-dot formats such as "26.441.74"); for the cells whose new
# text would otherwise be auto-parsed as a real number, force Text entry by
# flipping the NumberFormat to "@" first, then clear the formatting again so
# the cell keeps its original (unstyled) appearance while the stored value
# remains literal text.

$ws.Range("D2").Value = "26.441.74"
$ws.Range("E2").Value = "  +1.49%  "

$ws.Range("D3").Value = "1.693.74"
$ws.Range("E3").Value = "  +1.57%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.65%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.12"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5548"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +8.79%  "

$ws.Range("E7").Value = "  +0.60%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2719"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.68%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06489"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.61%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.17"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07644"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.77%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.563"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.18%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5828"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.43%  "

$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.000008472"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.17%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.332.65"
$ws.Range("E15").Value = "  -20.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.30"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.87%  "

$ws.Range("D17").Value = "26.483.12"
$ws.Range("E17").Value = "  +2.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.965"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.011"
$ws.Range("D19").ClearFormats()

$ws.Range("E20").Value = "  +1.90%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.53"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.261"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.15%  "

$ws.Range("E23").Value = "  +0.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "150.47"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1310"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +7.97%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.911"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.91%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.79"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.83%  "

$ws.Range("E28").Value = "  +7.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06338"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.13%  "

$ws.Range("E30").Value = "  +1.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.598"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.599"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.680"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.23%  "

$ws.Range("E34").Value = "  +2.90%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6227"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.93%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.406"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.61%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.724"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.236"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.55%  "

$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.123.64"
$ws.Range("E39").Value = "  +2.34%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01645"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8854"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.96%  "

$ws.Range("E42").Value = "  +0.72%  "

$ws.Range("E43").Value = "  -0.72%  "

$ws.Range("D44").Value = "1.844.17"
$ws.Range("E44").Value = "  +1.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000111"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.48%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.61"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.229"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.007"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05286"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4302"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.081"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.57%  "

